$wb = $excel.ActiveWorkbook
$taxon = $wb.Worksheets.Item("Taxon")
$ws = $wb.Worksheets.Add($null, $taxon)
$ws.Name = "Environment"

$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "env"

$ws.Range("A2").Value = "Name"

$ws.Range("A3").Value = "Temperature"
$ws.Range("B3").Value = 37

$ws.Range("A4").Value = "Temperature units"
$ws.Range("B4").Value = "C"

$ws.Range("A5").Value = "pH"
$ws.Range("B5").Value = 7.75

$ws.Range("A6").Value = "Database references"
$ws.Range("A7").Value = "Comments"
$ws.Range("A8").Value = "References"

$ws.Activate()
